$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formato PF")

# Row 3
$ws.Range("AW3").Value = 12345.22
$ws.Range("BK3").Value = 1234.5
$ws.Range("BR3").Value = 1234.7
$ws.Range("BS3").Value = 1234567.4
$ws.Range("BT3").Value = 12345.24
$ws.Range("CH3").Value = 12345.3
$ws.Range("CL3").Value = 1234567.7

# Row 4
$ws.Range("AW4").Value = 123456.8
$ws.Range("BK4").Value = 567.66
$ws.Range("BR4").Value = 5678.78
$ws.Range("BU4").Value = 75.2
$ws.Range("CE4").Value = 123.2
$ws.Range("CF4").Value = 12345.7

# Row 5
$ws.Range("BK5").Value = 8910.549999999999
